# gelaendeList.xlsx edit
#
# Summary of the change (see commit "Add files via upload"):
#   - Column B ("close" date) values 3-7 gain a fractional "end of day"
#     component (23:59:59 -> .9999884259 of a day) instead of plain
#     midnight-dated serials.
#   - The shared date number format (numFmtId 165, used by both the A
#     and B "open"/"close" columns on rows 3-7) is changed from the
#     Japanese "M月D日" pattern to an ISO "YYYY-MM-DD" pattern, which is
#     why both columns end up re-rendered even though only column B's
#     underlying values changed.
#   - The sheet's recorded dimension/used-range grows to A1:D1000 and the
#     remembered cursor/selection moves to B7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "close" date values (column B, rows 3-7) -------------------
# Same calendar day as before, but now carrying the end-of-day fractional
# time component (44647 -> 44647.9999884259, etc.)
$ws.Range("B3").Value2 = 44647.9999884259
$ws.Range("B4").Value2 = 44661.9999884259
$ws.Range("B5").Value2 = 44668.9999884259
$ws.Range("B6").Value2 = 44661.9999884259
$ws.Range("B7").Value2 = 44689.9999884259

# --- Re-format the date columns --------------------------------------------
# A3:B7 all share the same style/number format in the original workbook, so
# switching that format to an ISO date pattern touches both columns.
$ws.Range("A3:B7").NumberFormat = "YYYY\-MM\-DD"

# --- Grow the sheet's used range out to row 1000 ----------------------------
# Touching a cell far below the data (without changing its value) nudges the
# worksheet's recorded dimension to A1:D1000, matching the new extent, while
# leaving the existing data/styles untouched.
$ws.Range("D1000").NumberFormat = $ws.Range("D1000").NumberFormat

# --- Update the remembered selection ----------------------------------------
$ws.Range("B7").Select()
